# "add hole depth data, R2 B3 bottle weights"
#
# The "Rainfall 2" sheet tracks, per sample bottle, vol_water_applied (D),
# bottle_mass (E), water_plus_bottle_mass_collected (F) and a computed
# water_mass = F - E (G, shared formula). Rows 38-55 (samples C37R2..C54R2)
# were missing their bottle_mass(g) reading in column E; fill them in with
# the newly-recorded values. Column G recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rainfall 2")

$boxWeights = [ordered]@{
    38 = 88
    39 = 89
    40 = 88
    41 = 89
    42 = 104
    43 = 89
    44 = 90
    45 = 90
    46 = 90
    47 = 88
    48 = 101
    49 = 90
    50 = 89
    51 = 89
    52 = 89
    53 = 88
    54 = 88
    55 = 89
}

foreach ($row in $boxWeights.Keys) {
    $ws.Cells.Item($row, 5).Value = $boxWeights[$row]
}

# Match the author's ending cursor position/scroll in the saved file.
$ws.Activate()
$ws.Range("E56").Select()
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
